# onInterpolation 함수 & forceChangeResolution
# Update the long-text localization example strings on the "ui" sheet
# and move the active selection from C5 to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ui")

# Refresh the example strings used to explain the A/S key interpolation shortcut.
$ws.Range("B4").Value = "Press A to make Korean, press S to make English"
$ws.Range("C4").Value = "A를 눌러서 한국어로 바꾸고, S키눌러서 영어로 바꿔보세요."

# Make the sheet active and move the selected cell to B5.
$ws.Activate()
$ws.Range("B5").Select()
